$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 116670070
$ws.Cells.Item(70, 9).Value = 125001500
$ws.Cells.Item(70, 11).Value = 375004500
$ws.Cells.Item(70, 13).Value = -375004230
$ws.Cells.Item(73, 8).Value = 116670070
$ws.Cells.Item(73, 9).Value = 125001500
$ws.Cells.Item(73, 11).Value = 375004500
$ws.Cells.Item(73, 13).Value = -375003564
$ws.Cells.Item(86, 8).Value = 64816840
$ws.Cells.Item(86, 9).Value = 74076190
$ws.Cells.Item(86, 10).Value = 18520086
$ws.Cells.Item(86, 11).Value = 74076190
$ws.Cells.Item(86, 12).Value = 18520086
$ws.Cells.Item(86, 13).Value = -74075067
$ws.Cells.Item(86, 14).Value = -18522332
$ws.Cells.Item(89, 8).Value = 64816840
$ws.Cells.Item(89, 9).Value = 74076190
$ws.Cells.Item(89, 10).Value = 18520086
$ws.Cells.Item(89, 11).Value = 370380950
$ws.Cells.Item(89, 12).Value = 92600430
$ws.Cells.Item(89, 13).Value = -370375334
$ws.Cells.Item(89, 14).Value = -92611662
$ws.Cells.Item(100, 8).Value = 2328.3076
$ws.Cells.Item(100, 9).Value = 1614.2858
$ws.Cells.Item(100, 10).Value = 3161.3333
$ws.Cells.Item(100, 11).Value = 1614.2858
$ws.Cells.Item(100, 12).Value = 3161.3333
$ws.Cells.Item(100, 13).Value = -1073.2858
$ws.Cells.Item(100, 14).Value = -4243.3333
$ws.Cells.Item(132, 8).Value = 1391.9259
$ws.Cells.Item(132, 9).Value = 1559.35
$ws.Cells.Item(132, 10).Value = 913.5714
$ws.Cells.Item(132, 11).Value = 4678.049999999999
$ws.Cells.Item(132, 12).Value = 2740.7142
$ws.Cells.Item(132, 13).Value = -2148.049999999999
$ws.Cells.Item(132, 14).Value = -7800.7142
$ws.Cells.Item(137, 8).Value = 2905.1304
$ws.Cells.Item(137, 9).Value = 3326.4443
$ws.Cells.Item(137, 11).Value = 9979.332900000001
$ws.Cells.Item(137, 13).Value = -7429.332900000001
$ws.Cells.Item(138, 8).Value = 5184.878
$ws.Cells.Item(138, 10).Value = 11116.529
$ws.Cells.Item(138, 12).Value = 33349.587
$ws.Cells.Item(138, 14).Value = -43629.587
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2781253.2
$ws.Cells.Item(32, 9).Value = 2979462.5
$ws.Cells.Item(32, 11).Value = 2979462.5
$ws.Cells.Item(32, 13).Value = -2979175.5
$ws.Cells.Item(63, 8).Value = 1825
$ws.Cells.Item(63, 9).Value = 1790
$ws.Cells.Item(63, 10).Value = 2000
$ws.Cells.Item(63, 11).Value = 1790
$ws.Cells.Item(63, 12).Value = 2000
$ws.Cells.Item(63, 13).Value = -1104
$ws.Cells.Item(63, 14).Value = -3372
$ws.Cells.Item(66, 8).Value = 1825
$ws.Cells.Item(66, 9).Value = 1790
$ws.Cells.Item(66, 10).Value = 2000
$ws.Cells.Item(66, 11).Value = 8950
$ws.Cells.Item(66, 12).Value = 2000
$ws.Cells.Item(66, 13).Value = -5518
$ws.Cells.Item(66, 14).Value = -16864
$ws.Cells.Item(74, 8).Value = 50661.12
$ws.Cells.Item(74, 9).Value = 62945.73
$ws.Cells.Item(74, 11).Value = 62945.73
$ws.Cells.Item(74, 13).Value = -62071.73
$ws.Cells.Item(77, 8).Value = 50661.12
$ws.Cells.Item(77, 9).Value = 62945.73
$ws.Cells.Item(77, 11).Value = 314728.65
$ws.Cells.Item(77, 13).Value = -310360.65
$ws.Cells.Item(132, 8).Value = 6092.0444
$ws.Cells.Item(132, 9).Value = 4423.0967
$ws.Cells.Item(132, 10).Value = 9787.571
$ws.Cells.Item(132, 11).Value = 13269.2901
$ws.Cells.Item(132, 12).Value = 29362.713
$ws.Cells.Item(132, 13).Value = -10739.2901
$ws.Cells.Item(132, 14).Value = -34422.713
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4789.579
$ws.Cells.Item(16, 9).Value = 3000.875
$ws.Cells.Item(16, 10).Value = 6090.4546
$ws.Cells.Item(16, 11).Value = 3000.875
$ws.Cells.Item(16, 12).Value = 6090.4546
$ws.Cells.Item(16, 13).Value = -2713.875
$ws.Cells.Item(16, 14).Value = -6664.4546
$ws.Cells.Item(31, 8).Value = 7632.6885
$ws.Cells.Item(31, 9).Value = 3421.3103
$ws.Cells.Item(31, 11).Value = 3421.3103
$ws.Cells.Item(31, 13).Value = -3126.3103
$ws.Cells.Item(34, 8).Value = 7632.6885
$ws.Cells.Item(34, 9).Value = 3421.3103
$ws.Cells.Item(34, 11).Value = 3421.3103
$ws.Cells.Item(34, 13).Value = -3219.3103
$ws.Cells.Item(113, 8).Value = 4789.579
$ws.Cells.Item(113, 9).Value = 3000.875
$ws.Cells.Item(113, 10).Value = 6090.4546
$ws.Cells.Item(113, 11).Value = 3000.875
$ws.Cells.Item(113, 12).Value = 6090.4546
$ws.Cells.Item(113, 13).Value = -830.875
$ws.Cells.Item(113, 14).Value = -10430.4546
$ws.Cells.Item(134, 8).Value = 23818642
$ws.Cells.Item(134, 9).Value = 11487.134
$ws.Cells.Item(134, 10).Value = 37044840
$ws.Cells.Item(134, 11).Value = 34461.402
$ws.Cells.Item(134, 12).Value = 111134520
$ws.Cells.Item(134, 13).Value = -31926.402
$ws.Cells.Item(134, 14).Value = -111139590
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 201055.5
$ws.Cells.Item(137, 9).Value = 143665.28
$ws.Cells.Item(137, 11).Value = 430995.84
$ws.Cells.Item(137, 13).Value = -425895.84
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 112885.555
$ws.Cells.Item(80, 10).Value = 252045
$ws.Cells.Item(80, 12).Value = 252045
$ws.Cells.Item(80, 14).Value = -254041
$ws.Cells.Item(83, 8).Value = 112885.555
$ws.Cells.Item(83, 10).Value = 252045
$ws.Cells.Item(83, 12).Value = 1260225
$ws.Cells.Item(83, 14).Value = -1270209
$ws.Cells.Item(113, 8).Value = 6676.5713
$ws.Cells.Item(113, 10).Value = 7402.091
$ws.Cells.Item(113, 12).Value = 7402.091
$ws.Cells.Item(113, 14).Value = -11742.091
$ws.Cells.Item(126, 8).Value = 2908.9443
$ws.Cells.Item(126, 9).Value = 2859.5
$ws.Cells.Item(126, 11).Value = 8578.5
$ws.Cells.Item(126, 13).Value = -6108.5
$ws.Cells.Item(132, 8).Value = 1854
$ws.Cells.Item(132, 9).Value = 1826.5116
$ws.Cells.Item(132, 10).Value = 1985.3334
$ws.Cells.Item(132, 11).Value = 5479.5348
$ws.Cells.Item(132, 12).Value = 5956.0002
$ws.Cells.Item(132, 13).Value = -2949.5348
$ws.Cells.Item(132, 14).Value = -11016.0002
$ws.Cells.Item(134, 8).Value = 99997.664
$ws.Cells.Item(134, 10).Value = 99997.664
$ws.Cells.Item(134, 12).Value = 299992.992
$ws.Cells.Item(134, 14).Value = -305062.992
$ws.Cells.Item(136, 8).Value = 67032.60000000001
$ws.Cells.Item(136, 10).Value = 91721
$ws.Cells.Item(136, 12).Value = 275163
$ws.Cells.Item(136, 14).Value = -280263
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7900.2666
$ws.Cells.Item(7, 9).Value = 6299.5
$ws.Cells.Item(7, 11).Value = 6299.5
$ws.Cells.Item(7, 13).Value = -6187.5
$ws.Cells.Item(40, 8).Value = 5649.514
$ws.Cells.Item(40, 9).Value = 4951.115
$ws.Cells.Item(40, 10).Value = 7667.1113
$ws.Cells.Item(40, 11).Value = 4951.115
$ws.Cells.Item(40, 12).Value = 7667.1113
$ws.Cells.Item(40, 13).Value = -4815.115
$ws.Cells.Item(40, 14).Value = -7939.1113
$ws.Cells.Item(46, 8).Value = 12348039
$ws.Cells.Item(46, 10).Value = 12348039
$ws.Cells.Item(46, 12).Value = 12348039
$ws.Cells.Item(46, 14).Value = -12348415
$ws.Cells.Item(68, 8).Value = 2979
$ws.Cells.Item(68, 9).Value = 2969
$ws.Cells.Item(68, 11).Value = 2969
$ws.Cells.Item(68, 13).Value = -2220
$ws.Cells.Item(71, 8).Value = 2979
$ws.Cells.Item(71, 9).Value = 2969
$ws.Cells.Item(71, 11).Value = 14845
$ws.Cells.Item(71, 13).Value = -11101
$ws.Cells.Item(82, 8).Value = 1831.1
$ws.Cells.Item(82, 9).Value = 1643.2
$ws.Cells.Item(82, 10).Value = 2019
$ws.Cells.Item(82, 11).Value = 1643.2
$ws.Cells.Item(82, 12).Value = 2019
$ws.Cells.Item(82, 13).Value = -1282.2
$ws.Cells.Item(82, 14).Value = -2741
$ws.Cells.Item(85, 8).Value = 1831.1
$ws.Cells.Item(85, 9).Value = 1643.2
$ws.Cells.Item(85, 10).Value = 2019
$ws.Cells.Item(85, 11).Value = 1643.2
$ws.Cells.Item(85, 12).Value = 2019
$ws.Cells.Item(85, 13).Value = -395.2
$ws.Cells.Item(85, 14).Value = -4515
$ws.Cells.Item(98, 8).Value = 54796
$ws.Cells.Item(98, 10).Value = 54796
$ws.Cells.Item(98, 12).Value = 54796
$ws.Cells.Item(98, 14).Value = -60786
$ws.Cells.Item(126, 8).Value = 7900.2666
$ws.Cells.Item(126, 9).Value = 6299.5
$ws.Cells.Item(126, 11).Value = 18898.5
$ws.Cells.Item(126, 13).Value = -16428.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 673.875
$ws.Cells.Item(100, 9).Value = 673.875
$ws.Cells.Item(100, 11).Value = 1347.75
$ws.Cells.Item(100, 13).Value = -806.75
$ws.Cells.Item(122, 8).Value = 4900.231
$ws.Cells.Item(122, 9).Value = 3025.5
$ws.Cells.Item(122, 10).Value = 7899.8
$ws.Cells.Item(122, 11).Value = 9076.5
$ws.Cells.Item(122, 12).Value = 23699.4
$ws.Cells.Item(122, 13).Value = -6626.5
$ws.Cells.Item(122, 14).Value = -28599.4
